$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "60.678.97"
$ws.Range("E2").Value = "  +4.16%  "
$ws.Range("D3").Value = "2.628.67"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.06"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +6.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.92"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.606"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +4.11%  "
$ws.Range("D9").Value = "2.649.17"
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.84"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("E11").Value = "  +5.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.150"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +8.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.345"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +4.39%  "
$ws.Range("D14").Value = "3.093.23"
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("D15").Value = "60.603.83"
$ws.Range("E15").Value = "  +4.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.11"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +7.23%  "
$ws.Range("E17").Value = "  +3.45%  "
$ws.Range("D18").Value = "2.630.81"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.91"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.42"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.37"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +3.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.19"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.439"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +4.83%  "
$ws.Range("E26").Value = "  +4.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("E28").Value = "  +4.75%  "
$ws.Range("D29").Value = "0.0₃0803"
$ws.Range("E29").Value = "  +10.44%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.71"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +4.83%  "
$ws.Range("E32").Value = "  +4.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.33"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.15"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.10"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +5.56%  "
$ws.Range("E36").Value = "  +5.50%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.887"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +8.61%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.880"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +6.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.52"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.52"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +7.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "300.04"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +5.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.65"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  +4.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.603"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0544"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.46"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +5.67%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.70"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.28"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +15.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0235"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +4.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.67"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +6.76%  "
